$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: Clean up the "14 What is a primary key?" answer paragraph.
#   Remove the gramStart/gramEnd proofing markers around "The" and merge
#   that run together with the trailing " primary key ... null" run into a
#   single run, while leaving the preceding " " run intact.
# ---------------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute( `
    "identifier. The primary key of a relational database must be unique. Every row of data must have a primary key value and none of the rows can be null", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "identifier. The primary key of a relational database must be unique. Every row of data must have a primary key value and none of the rows can be null", `
    2)

$answerPara = $d.Paragraphs(87)
$paraStart = $answerPara.Range.Start
$paraText = $answerPara.Range.Text
$spaceOffset = $paraStart + $paraText.IndexOf(" The primary key")
$theOffset = $spaceOffset + 1
$paraEnd = $answerPara.Range.End - 1

# Re-split the run boundary right before "The" so the leading " " run stays
# separate from the newly merged "The primary key ... null" run.
$afterThe = $d.Range($theOffset, $paraEnd)
$afterThe.Font.Bold = $true
$afterThe.Font.Bold = $false

$spaceRun = $d.Range($spaceOffset, $theOffset)
$spaceRun.Font.Bold = $true
$spaceRun.Font.Bold = $false

# ---------------------------------------------------------------------------
# Step 2: Insert the new Q15 heading paragraph and its answer paragraph right
# after the "14" answer paragraph (and before the trailing blank paragraph).
# ---------------------------------------------------------------------------
$answerPara = $d.Paragraphs(87)
$answerPara.Range.InsertParagraphAfter()
$headingPara = $d.Paragraphs(88)
$headingPara.Range.InsertParagraphAfter()

# Fill in the plain answer paragraph first, while it is still unformatted, so
# it does not inherit the bold/large heading formatting.
$bodyPara = $d.Paragraphs(89)
$bodyPara.Range.Text = "In SQL, the GROUP BY clause is used to create summary rows out of rows that have the same values in a set of specified columns. In order to do computations on groups of rows as opposed to individual rows, it is frequently used in conjunction with aggregate functions like SUM, COUNT, AVG, MAX, or MIN. we may produce summary reports and perform more in-depth data analysis using the GROUP BY clause."

# Fill in the bold heading paragraph: "15" + "." + " What is the purpose of
# the GROUP BY clause, and how is it used? " as three separate runs.
$headingPara = $d.Paragraphs(88)
$headingRange = $headingPara.Range
$headingRange.Text = "15"
$headingRange.Font.Bold = $true
$headingRange.Font.Size = 14

$afterFifteen = $d.Paragraphs(88).Range
$dotInsert = $d.Range($afterFifteen.End - 1, $afterFifteen.End - 1)
$dotInsert.InsertAfter(".")

$afterDot = $d.Paragraphs(88).Range
$restInsert = $d.Range($afterDot.End - 1, $afterDot.End - 1)
$restInsert.InsertAfter(" What is the purpose of the GROUP BY clause, and how is it used? ")

$headingRangeFinal = $d.Paragraphs(88).Range
$headingStart = $headingRangeFinal.Start
$splitPoint = $headingStart + 2
$splitRun = $d.Range($splitPoint, $splitPoint + 1)
$splitRun.Font.Bold = $false
$splitRunBack = $d.Range($splitPoint, $splitPoint + 1)
$splitRunBack.Font.Bold = $true

Write-Output "Paragraph 87: $($d.Paragraphs(87).Range.Text)"
Write-Output "Paragraph 88: $($d.Paragraphs(88).Range.Text)"
Write-Output "Paragraph 89: $($d.Paragraphs(89).Range.Text)"
Write-Output "Total paragraphs: $($d.Paragraphs.Count)"
